# Create a EvidenceDoc oriented to objects
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status/result values
$ws.Range("B2").Value = "No"
$ws.Range("C2").Value = ""
$ws.Range("H2").Value = "25/05/2020"

$ws.Range("B3").Value = "Yes"
$ws.Range("H3").Value = "25/05/2020"

$ws.Range("B4").Value = "No"
$ws.Range("C4").Value = ""
$ws.Range("H4").Value = "25/05/2020"

$ws.Range("H5").Value = "28_04_2020--21_22_15 376"

$ws.Range("B6").Value = "No"
$ws.Range("C6").Value = ""
$ws.Range("H6").Value = "25/05/2020"

# Move the active cell selection
$ws.Range("B10").Select()
